$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.082.38"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.891.46"
$ws.Range("E3").Value = "  +1.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9984"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.53"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5143"
$ws.Range("E7").Value = "  +1.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3753"
$ws.Range("E8").Value = "  +3.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07204"
$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.21"
$ws.Range("E10").Value = "  +2.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9046"
$ws.Range("E11").Value = "  +1.14%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07644"
$ws.Range("E12").Value = "  +2.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.894.81"
$ws.Range("E13").Value = "  +1.80%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.86"
$ws.Range("E14").Value = "  +2.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.263"
$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9989"
$ws.Range("E16").Value = "  -0.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008486"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.43"
$ws.Range("E18").Value = "  +2.09%  "

$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.107.32"
$ws.Range("E20").Value = "  +0.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.068"
$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.124.23"
$ws.Range("E22").Value = "  +1.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.56"
$ws.Range("E23").Value = "  +1.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.402"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("E25").Value = "  +10.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.67"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.773"
$ws.Range("E27").Value = "  -1.04%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.05"
$ws.Range("E28").Value = "  +0.95%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.42"
$ws.Range("E29").Value = "  +1.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.949"
$ws.Range("E30").Value = "  +5.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.828"
$ws.Range("E31").Value = "  +3.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09175"
$ws.Range("E32").Value = "  -0.99%  "

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.237"
$ws.Range("E34").Value = "  +7.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7828"
$ws.Range("E35").Value = "  +4.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.986"
$ws.Range("E36").Value = "  +0.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.288"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.613"
$ws.Range("E38").Value = "  +4.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01997"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5596"
$ws.Range("E40").Value = "  +1.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075"
$ws.Range("E41").Value = "  +0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.089"
$ws.Range("E42").Value = "  +7.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.649"
$ws.Range("E43").Value = "  +2.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "117.86"
$ws.Range("E44").Value = "  -0.54%  "

$ws.Range("E45").Value = "  +2.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4801"
$ws.Range("E46").Value = "  +2.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.25"
$ws.Range("E47").Value = "  +2.30%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9978"
$ws.Range("E48").Value = "  -0.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.600"
$ws.Range("E49").Value = "  +2.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.52"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.08"
$ws.Range("E51").Value = "  +1.50%  "
